$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns (in order) that receive new per-row values; column A, G and H are untouched.
$cols = @("B","C","D","E","F","I","J","K","L","M","N")

# Map of row number -> new values for @cols, taken from the updated 380kV case run.
$rowData = @{
    2 = @(1.02, 1.004727430580851, 1.026513451475777, 1.007762213487518, 1.002867825820396, 1.028140619206577, 1.010014891386621, 1.029335878222202, 1.010640516912302, 1.005761089072324, 1.007465362557762)
    3 = @(1.02, 1.005822002272901, 1.026832767717716, 1.008694518224989, 1.004617363073518, 1.028119731977348, 1.010738946473879, 1.029464047931005, 1.01137594559515, 1.007310333110168, 1.007712251451404)
    4 = @(1.02, 1.006530092566826, 1.027038990663314, 1.009297975763954, 1.005748917304874, 1.02810454699596, 1.01120681070201, 1.029545830079758, 1.011851397160795, 1.008311870664742, 1.007871580600961)
    5 = @(1.02, 1.006827735276856, 1.027125589626558, 1.009551717083347, 1.006224509550559, 1.028097762610936, 1.011403347421962, 1.029579933844576, 1.012051177842443, 1.008732703681971, 1.007938461346076)
    6 = @(1.02, 1.006877708613519, 1.027140124207022, 1.00959432420678, 1.006304357228488, 1.028096599963333, 1.011436337855459, 1.029585643700015, 1.012084716123469, 1.008803351132689, 1.007949684992606)
    7 = @(1.02, 1.006534069831341, 1.027040148187686, 1.009301366078558, 1.005755272624287, 1.028104457918005, 1.01120943743596, 1.029546286868166, 1.011854067027955, 1.008317494683953, 1.007872474662037)
    8 = @(1.02, 1.005097383039582, 1.026621446480098, 1.008077250859506, 1.003459202297312, 1.02813390540539, 1.010259724350782, 1.029379431388346, 1.010889146272544, 1.006284860424486, 1.00754888794167)
    9 = @(1.02, 1.002564339795591, 1.025880711243438, 1.005921629741186, 0.9994088644037524, 1.028173049635261, 1.008581170752679, 1.029076656696341, 1.009185553564632, 1.002695617146474, 1.00697541805312)
    10 = @(1.02, 1.000874535888566, 1.025385058414673, 1.004485421829613, 0.9967050907846262, 1.028190631074268, 1.007458639509106, 1.02886901662911, 1.008047530887223, 1.000297222101511, 1.006590881253872)
    11 = @(1.02, 1.000142534010965, 1.025170031100778, 1.003863715011103, 0.995533349279506, 1.028196236062541, 1.006971717307854, 1.028777753015332, 1.007554189785336, 0.999257251006217, 1.006423839134035)
    12 = @(1.02, 0.99987058688519, 1.02509010175386, 1.003632810757945, 0.9950979518761309, 1.028198017265132, 1.006790721502815, 1.028743651841878, 1.007370853485298, 0.9988707312407346, 1.006361711295728)
    13 = @(1.02, 0.9999289227360615, 1.025107249484421, 1.00368233938652, 0.9951913535466643, 1.028197648785465, 1.006829551690842, 1.028750975765278, 1.007410183731069, 0.9989536515352102, 1.006375041600458)
    14 = @(1.02, 1.000120055795136, 1.025163425298077, 1.003844627887068, 0.9954973625354055, 1.028196389426981, 1.006956758827764, 1.028774938309027, 1.007539036944252, 0.9992253058768764, 1.006418705282296)
    15 = @(1.02, 1.000237812589679, 1.025198029374264, 1.003944622492043, 0.9956858831607941, 1.028195573671256, 1.007035117889274, 1.028789675730796, 1.007618416018253, 0.9993926505345242, 1.006445597173323)
    16 = @(1.02, 1.000923109737314, 1.025399320714895, 1.004526686116763, 0.9967828335567679, 1.028190216864643, 1.007490936651406, 1.028875045105079, 1.008080260149639, 1.000366210124534, 1.006601955966833)
    17 = @(1.02, 1.001352894509748, 1.025525478374841, 1.00489184635577, 0.997470648322167, 1.028186319765488, 1.007776628254542, 1.02892823368536, 1.00836980892789, 1.000976502807035, 1.006699892041016)
    18 = @(1.02, 1.00160355141928, 1.025599024633332, 1.005104855702431, 0.9978717443451044, 1.02818385285959, 1.007943184624868, 1.02895912692504, 1.008538642982987, 1.001332336979295, 1.006756964900417)
    19 = @(1.02, 1.001689014075675, 1.025624095254975, 1.005177489416966, 0.9980084920162302, 1.028182978809425, 1.007999962080195, 1.028969638492465, 1.008596201784786, 1.001453643932613, 1.006776416535688)
    20 = @(1.02, 1.001306785723031, 1.025511946906694, 1.004852666315862, 0.9973968621890754, 1.028186757925733, 1.007745984820907, 1.028922540567382, 1.008338748771612, 1.000911038655298, 1.006689389768212)
    21 = @(1.02, 1.000063773251871, 1.025146884520152, 1.003796837283646, 0.9954072550388525, 1.028196768571164, 1.006919303136536, 1.028767887494854, 1.007501095313476, 0.9991453168265775, 1.006405849655614)
    22 = @(1.02, 0.9992819567113377, 1.024917017289605, 1.003133141372291, 0.9941553759467042, 1.028201323127522, 1.006398775255053, 1.028669484195791, 1.006973922460185, 0.9980338120025594, 1.006227108057501)
    23 = @(1.02, 0.9996964405381674, 1.025038905493229, 1.003484965928803, 0.9948191132268769, 1.028199073251911, 1.006674789740152, 1.028721759700743, 1.007253435431472, 0.9986231708395156, 1.006321906973569)
    24 = @(1.02, 1.001327620389373, 1.025518061316702, 1.004870370035164, 0.9974302032580508, 1.028186560538906, 1.007759831525004, 1.028925113446604, 1.008352783688695, 1.000940619516847, 1.006694135452648)
    25 = @(1.02, 1.003219376118636, 1.026072541693807, 1.006478749260445, 1.000456554253951, 1.028164434695646, 1.009015725207216, 1.029155958733625, 1.009626370872034, 1.0036244642001, 1.007124063766357)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $colLetter = $cols[$i]
        $ws.Range("$colLetter$row").Value2 = $values[$i]
    }
}
